$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new plain_password value (column C), rows 2 through 23
$newValues = @{
    2  = "LQr975"
    3  = "D4nA89"
    4  = "rFnc57"
    5  = "Adjw92"
    6  = "GbPn87"
    7  = "B4HF61"
    8  = "onn418"
    9  = "kHcD56"
    10 = "Ph0C90"
    11 = "TUR169"
    12 = "Dc5A19"
    13 = "NJi027"
    14 = "Fvsy80"
    15 = "2arf25"
    16 = "I9xV63"
    17 = "wSUr83"
    18 = "wpmm46"
    19 = "PZXC80"
    20 = "GVlt95"
    21 = "KFqv86"
    22 = "cYhB18"
    23 = "FpNf54"
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
